# Adds "Buoc 11" / "Buoc 12" / "Buoc 13" (new Laravel UI / npm steps)
# before the final "Buoc 11 (cuoi cung)" paragraph, and renumbers that
# paragraph's step number from "1"+"1" (split across two runs) into a
# single "13" run.

$d = $word.ActiveDocument

# Locate the target paragraph robustly (its old text is unique: it
# mentions "php artisan serve", which is plain ASCII).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "php artisan serve") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph (the final 'php artisan serve' step)"
}

$target = $d.Paragraphs.Item($targetIndex)
$full = $target.Range

# Replace the whole paragraph (content + its paragraph mark) with four
# paragraphs: the first keeps the original paragraph identity, the
# remaining three are brand-new paragraphs inserted ahead of the
# renumbered "(cuoi cung)" step.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="03FEDCB0" w14:textId="173C6149" w:rsidR="009B0FD4" w:rsidRPr="009B0FD4" w:rsidRDefault="009B0FD4" w:rsidP="009B0FD4"><w:pPr><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r w:rsidRPr="009B0FD4"><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">Bước 11: Nhập: </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t>composer require laravel/ui --dev</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">Bước 12: Nhập: </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">php artisan ui vue </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t>auth</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>Bước</w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> 13: </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">Nhập </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t>npm install &amp;&amp; npm run dev</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">Bước </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t>13</w:t></w:r><w:r w:rsidRPr="009B0FD4"><w:rPr><w:noProof/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> (cuối cùng): Nhập: php artisan serve rồi </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$full.InsertXML($xml)

Write-Host "Inserted steps 11-13 and renumbered final step at paragraph $targetIndex"
